$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.687.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.506.47'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.61'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.01'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.506.15'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.64'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.403'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.098.93'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000200'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.64'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.541.16'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.656.62'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.91'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -5.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '412.14'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.594'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.33'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.644.93'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000115'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.68'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.501.97'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.153'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.21'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.46'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '175.22'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.25'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -14.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.18'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -7.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.57'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -8.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0815'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.02'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.851'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.23'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -8.17%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.39'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.06'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.21'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.08'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -9.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.75'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -9.48%  '
